# gradeList.xlsx: clear the existing numeric grade for student sv23123
# (B2) so the board reflects the updated/total grade workflow.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = ""
